# fix: fk_id of tambon data
#
# The id/name lookup table that lives in A1:B7 (header "id"/"name" in row 1,
# then 6 region rows) gets duplicated directly below itself, in A8:B14, so
# the sheet grows from A1:B7 to A1:B14.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$src = $ws.Range("A1:B7")
$dst = $ws.Range("A8:B14")

$src.Copy()
$dst.PasteSpecial(-4163)  # xlPasteValues

$excel.CutCopyMode = $false
